$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values - B2 and D2 are cleared (deleted), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -12.049607504453629
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -11.499985915910047

# Row 3 values updated
$ws.Range("B3").Value = -16.827364875250563
$ws.Range("C3").Value = -7.3081369995935033
$ws.Range("D3").Value = -19.35657950153503
$ws.Range("E3").Value = 17.973387780569752

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
